# "Error Calculations and Plots"
#
# The sheet is a missing-data table (ID, A, B, C, D, F). This edit:
#   1. Removes two data rows entirely (ID "RM 232" and ID "SC 5"'s old
#      neighbour "SC 92"), which shifts every row below them up and
#      shrinks the used range from A1:F35 to A1:F33.
#   2. Re-imputes/clears a number of individual cells elsewhere in the
#      table (some values that were blank now have numbers filled in;
#      some values that had numbers are now blanked out).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: delete the two rows ------------------------------------------
# Row 26 is "RM 232". Deleting it shifts "SC 5" (old row 27) up to row 26,
# "SC 92" (old row 28) up to row 27, etc.
$ws.Rows(26).Delete()
# After the first delete, the old "SC 92" row is now row 27 - delete that too.
$ws.Rows(27).Delete()

# --- Step 2: per-cell value changes (using the post-deletion row numbers,
#     spreadsheet column letters A-F) ----------------------------------------

# RM 9 (row 4): clear E
$ws.Range("E4").Value = ""

# RM 14 (row 5): clear D
$ws.Range("D5").Value = ""

# RM 21 (row 6): fill in C
$ws.Range("C6").Value = 15.1

# RM 32 (row 7): clear E
$ws.Range("E7").Value = ""

# RM 38 (row 8): clear C
$ws.Range("C8").Value = ""

# RM 58 (row 11): fill in D
$ws.Range("D11").Value = -15.5

# RM 116 (row 17): fill in E
$ws.Range("E17").Value = -7.3

# RM 125 (row 19): fill in C, clear D
$ws.Range("C19").Value = 13.2
$ws.Range("D19").Value = ""

# RM 135 (row 21): clear C
$ws.Range("C21").Value = ""

# RM 140 (row 23): fill in C and D
$ws.Range("C23").Value = 12.2
$ws.Range("D23").Value = -13.9

# RM 142a (row 24): fill in E
$ws.Range("E24").Value = -8.1

# RM 145 (row 25): fill in D
$ws.Range("D25").Value = -15.5

# SC 5 (row 26, formerly RM 232's data slot + SC 5 label): clear B
$ws.Range("B26").Value = ""

# SC 101 (row 27): fill B, clear C and D, fill E
$ws.Range("B27").Value = -20.4
$ws.Range("C27").Value = ""
$ws.Range("D27").Value = ""
$ws.Range("E27").Value = -10

# SC 105 (row 28): fill E
$ws.Range("E28").Value = -5.9

# SC 119 (row 29): clear B, fill C, clear D and E
$ws.Range("B29").Value = ""
$ws.Range("C29").Value = 11.2
$ws.Range("D29").Value = ""
$ws.Range("E29").Value = ""

# SC 120 (row 30): fill D, clear E
$ws.Range("D30").Value = -13.6
$ws.Range("E30").Value = ""

# SC 193 (row 32): clear E
$ws.Range("E32").Value = ""

# SC 232 (row 33): fill D
$ws.Range("D33").Value = -14.1
